$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 196
$ws.Cells.Item(2, 2).Value = 'Blockaid'
$ws.Cells.Item(2, 3).Value = 'Enterprise Account Executive (Fintech)'
$ws.Cells.Item(2, 4).Value = 'Dennis Coombs'
$ws.Cells.Item(2, 5).Value = '2nd Interview'
$ws.Cells.Item(2, 6).Value = 45992
$ws.Cells.Item(2, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(3, 1).Value = 196
$ws.Cells.Item(3, 2).Value = 'Blockaid'
$ws.Cells.Item(3, 3).Value = 'Enterprise Account Executive (Fintech)'
$ws.Cells.Item(3, 4).Value = 'Erik Gallant'
$ws.Cells.Item(3, 5).Value = 'CV Sent'
$ws.Cells.Item(3, 6).Value = 45987
$ws.Cells.Item(3, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(4, 1).Value = 196
$ws.Cells.Item(4, 2).Value = 'Blockaid'
$ws.Cells.Item(4, 3).Value = 'Enterprise Account Executive (Fintech)'
$ws.Cells.Item(4, 4).Value = 'Madyson Almeida'
$ws.Cells.Item(4, 5).Value = '2nd Interview'
$ws.Cells.Item(4, 6).Value = 45994
$ws.Cells.Item(4, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(5, 1).Value = 196
$ws.Cells.Item(5, 2).Value = 'Blockaid'
$ws.Cells.Item(5, 3).Value = 'Enterprise Account Executive (Fintech)'
$ws.Cells.Item(5, 4).Value = 'Mike Gomez'
$ws.Cells.Item(5, 5).Value = '1st Interview'
$ws.Cells.Item(5, 6).Value = 45992
$ws.Cells.Item(5, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(6, 1).Value = 484
$ws.Cells.Item(6, 2).Value = 'Cognition AI'
$ws.Cells.Item(6, 3).Value = 'Enterprise Account Executive (US)'
$ws.Cells.Item(6, 4).Value = 'Tiffany Shih'
$ws.Cells.Item(6, 5).Value = '1st Interview'
$ws.Cells.Item(6, 6).Value = 45994
$ws.Cells.Item(6, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(7, 1).Value = 484
$ws.Cells.Item(7, 2).Value = 'Cognition AI'
$ws.Cells.Item(7, 3).Value = 'Enterprise Account Executive (US)'
$ws.Cells.Item(7, 4).Value = 'Mikaela Stamas'
$ws.Cells.Item(7, 5).Value = '2nd Interview'
$ws.Cells.Item(7, 6).Value = 45986
$ws.Cells.Item(7, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(8, 1).Value = 484
$ws.Cells.Item(8, 2).Value = 'Cognition AI'
$ws.Cells.Item(8, 3).Value = 'Enterprise Account Executive (US)'
$ws.Cells.Item(8, 4).Value = 'Matt Bartley'
$ws.Cells.Item(8, 5).Value = '1st Interview'
$ws.Cells.Item(8, 6).Value = 45977
$ws.Cells.Item(8, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(9, 1).Value = 484
$ws.Cells.Item(9, 2).Value = 'Cognition AI'
$ws.Cells.Item(9, 3).Value = 'Enterprise Account Executive (US)'
$ws.Cells.Item(9, 4).Value = 'Erik Abbott'
$ws.Cells.Item(9, 5).Value = '1st Interview'
$ws.Cells.Item(9, 6).Value = 45963
$ws.Cells.Item(9, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(10, 1).Value = 633
$ws.Cells.Item(10, 2).Value = 'Factory'
$ws.Cells.Item(10, 3).Value = 'CS1 Factory - Enterprise AE'
$ws.Cells.Item(10, 4).Value = 'Tiffany Shih'
$ws.Cells.Item(10, 5).Value = 'CV Sent'
$ws.Cells.Item(10, 6).Value = 45981
$ws.Cells.Item(10, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(11, 1).Value = 681
$ws.Cells.Item(11, 2).Value = 'Metaview'
$ws.Cells.Item(11, 3).Value = 'Metaview - UK MM / Enterprise AE'
$ws.Cells.Item(11, 4).Value = 'Marina Shynkarenka'
$ws.Cells.Item(11, 5).Value = 'CV Sent'
$ws.Cells.Item(11, 6).Value = 45975
$ws.Cells.Item(11, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(12, 1).Value = 696
$ws.Cells.Item(12, 2).Value = 'Cognition AI'
$ws.Cells.Item(12, 3).Value = 'Founding EMEA AE'
$ws.Cells.Item(12, 4).Value = 'Katie Pope'
$ws.Cells.Item(12, 5).Value = '1st Interview'
$ws.Cells.Item(12, 6).Value = 45993
$ws.Cells.Item(12, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(13, 1).Value = 750
$ws.Cells.Item(13, 2).Value = 'Novee.io'
$ws.Cells.Item(13, 3).Value = 'Enterprise Account Executive (US)'
$ws.Cells.Item(13, 4).Value = 'Chas Sheffield'
$ws.Cells.Item(13, 5).Value = 'CV Sent'
$ws.Cells.Item(13, 6).Value = 45986
$ws.Cells.Item(13, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(14, 1).Value = 750
$ws.Cells.Item(14, 2).Value = 'Novee.io'
$ws.Cells.Item(14, 3).Value = 'Enterprise Account Executive (US)'
$ws.Cells.Item(14, 4).Value = 'Rob Harvey'
$ws.Cells.Item(14, 5).Value = '3rd Interview'
$ws.Cells.Item(14, 6).Value = 45987
$ws.Cells.Item(14, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(15, 1).Value = 757
$ws.Cells.Item(15, 2).Value = 'Blockaid'
$ws.Cells.Item(15, 3).Value = 'SDR (Singapore)'
$ws.Cells.Item(15, 4).Value = 'Sean Hsien'
$ws.Cells.Item(15, 5).Value = 'CV Sent'
$ws.Cells.Item(15, 6).Value = 45992
$ws.Cells.Item(15, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(16, 1).Value = 757
$ws.Cells.Item(16, 2).Value = 'Blockaid'
$ws.Cells.Item(16, 3).Value = 'SDR (Singapore)'
$ws.Cells.Item(16, 4).Value = 'Dinie Mifdhal'
$ws.Cells.Item(16, 5).Value = '3rd Interview'
$ws.Cells.Item(16, 6).Value = 45985
$ws.Cells.Item(16, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(17, 1).Value = 762
$ws.Cells.Item(17, 2).Value = 'Energy Robotics'
$ws.Cells.Item(17, 3).Value = 'Account Executive (EMEA)'
$ws.Cells.Item(17, 4).Value = 'Shiwalla Singh'
$ws.Cells.Item(17, 5).Value = '2nd Interview'
$ws.Cells.Item(17, 6).Value = 45981
$ws.Cells.Item(17, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(18, 1).Value = 779
$ws.Cells.Item(18, 2).Value = 'Energy Robotics'
$ws.Cells.Item(18, 3).Value = 'SDR London'
$ws.Cells.Item(18, 4).Value = 'Toby Sproston'
$ws.Cells.Item(18, 5).Value = '1st Interview'
$ws.Cells.Item(18, 6).Value = 45981
$ws.Cells.Item(18, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(19, 1).Value = 779
$ws.Cells.Item(19, 2).Value = 'Energy Robotics'
$ws.Cells.Item(19, 3).Value = 'SDR London'
$ws.Cells.Item(19, 4).Value = 'Bilal Javaid'
$ws.Cells.Item(19, 5).Value = '2nd Interview'
$ws.Cells.Item(19, 6).Value = 45992
$ws.Cells.Item(19, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(20, 1).Value = 779
$ws.Cells.Item(20, 2).Value = 'Energy Robotics'
$ws.Cells.Item(20, 3).Value = 'SDR London'
$ws.Cells.Item(20, 4).Value = 'Daniel Murphy'
$ws.Cells.Item(20, 5).Value = '2nd Interview'
$ws.Cells.Item(20, 6).Value = 45982
$ws.Cells.Item(20, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(21, 1).Value = 779
$ws.Cells.Item(21, 2).Value = 'Energy Robotics'
$ws.Cells.Item(21, 3).Value = 'SDR London'
$ws.Cells.Item(21, 4).Value = 'Justin Kanapathy'
$ws.Cells.Item(21, 5).Value = '2nd Interview'
$ws.Cells.Item(21, 6).Value = 45992
$ws.Cells.Item(21, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(22, 1).Value = 779
$ws.Cells.Item(22, 2).Value = 'Energy Robotics'
$ws.Cells.Item(22, 3).Value = 'SDR London'
$ws.Cells.Item(22, 4).Value = 'Morgan Males'
$ws.Cells.Item(22, 5).Value = '3rd Interview'
$ws.Cells.Item(22, 6).Value = 45988
$ws.Cells.Item(22, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(23, 1).Value = 782
$ws.Cells.Item(23, 2).Value = 'Port'
$ws.Cells.Item(23, 3).Value = 'Enterprise AE'
$ws.Cells.Item(23, 4).Value = 'Jamie Summers'
$ws.Cells.Item(23, 5).Value = '1st Interview'
$ws.Cells.Item(23, 6).Value = 45981
$ws.Cells.Item(23, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(24, 1).Value = 783
$ws.Cells.Item(24, 2).Value = 'Port'
$ws.Cells.Item(24, 3).Value = 'Mid-Market AE'
$ws.Cells.Item(24, 4).Value = 'Shahz Shuja'
$ws.Cells.Item(24, 5).Value = '2nd Interview'
$ws.Cells.Item(24, 6).Value = 45994
$ws.Cells.Item(24, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(25, 1).Value = 813
$ws.Cells.Item(25, 2).Value = 'Laurel'
$ws.Cells.Item(25, 3).Value = 'Enterprise Account Executive UK x4'
$ws.Cells.Item(25, 4).Value = 'Luiz Kemmer'
$ws.Cells.Item(25, 5).Value = '1st Interview'
$ws.Cells.Item(25, 6).Value = 45992
$ws.Cells.Item(25, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(26, 1).Value = 816
$ws.Cells.Item(26, 2).Value = 'Allium'
$ws.Cells.Item(26, 3).Value = 'Enterprise Account Executive'
$ws.Cells.Item(26, 4).Value = 'Erik Hug'
$ws.Cells.Item(26, 5).Value = 'CV Sent'
$ws.Cells.Item(26, 6).Value = 45980
$ws.Cells.Item(26, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(27, 1).Value = 832
$ws.Cells.Item(27, 2).Value = 'Blockaid'
$ws.Cells.Item(27, 3).Value = 'Enterprise Account Executive'
$ws.Cells.Item(27, 4).Value = 'Rib Das'
$ws.Cells.Item(27, 5).Value = 'CV Sent'
$ws.Cells.Item(27, 6).Value = 45987
$ws.Cells.Item(27, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(28, 1).Value = 832
$ws.Cells.Item(28, 2).Value = 'Blockaid'
$ws.Cells.Item(28, 3).Value = 'Enterprise Account Executive'
$ws.Cells.Item(28, 4).Value = 'Jason Ong'
$ws.Cells.Item(28, 5).Value = 'CV Sent'
$ws.Cells.Item(28, 6).Value = 45987
$ws.Cells.Item(28, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(29, 1).Value = 832
$ws.Cells.Item(29, 2).Value = 'Blockaid'
$ws.Cells.Item(29, 3).Value = 'Enterprise Account Executive'
$ws.Cells.Item(29, 4).Value = 'Bonner Pang'
$ws.Cells.Item(29, 5).Value = 'CV Sent'
$ws.Cells.Item(29, 6).Value = 45987
$ws.Cells.Item(29, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
